$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.323.14"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.035.55"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "197.94"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "619.20"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("D10").Value = "3.035.80"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "5.28"
$ws.Range("E13").Value = "  +6.68%  "
$ws.Range("D14").Value = "3.593.71"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "28.72"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "76.194.33"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "0.0000192"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "3.040.87"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").Value = "13.45"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "8.91"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "378.06"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "4.33"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "3.193.86"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "72.85"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "4.32"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "9.70"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "8.24"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "490.22"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "20.52"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "162.16"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "20.03"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.117"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").Value = "190.86"
$ws.Range("E41").Value = "  +6.21%  "
$ws.Range("E42").Value = "  -5.11%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "0.797"
$ws.Range("E44").Value = "  +20.68%  "
$ws.Range("D45").Value = "5.05"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "41.89"
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "1.25"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "2.40"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "3.86"
$ws.Range("E51").Value = "  -2.22%  "
